$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the observation records that currently sit in rows 2, 3, 4
# and 6 (row 5 is untouched):
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 6   (row 6 also carried a few extra columns -
#                              J, K, L, N, AC, AF - that move along with it)
#   new row 6 <- old row 3   (and therefore no longer has those extra columns)

# Columns that carry real data and are present on every one of these rows.
$dataCols = @("A","B","C","D","E","F","G","H","P","Q","R","S","T","U","V","W","Y","Z","AA","AB","AD","AE","AG","AW","AX","AY")

# The extra columns that (before the edit) only exist on row 6. K and AC hold
# actual values; J, L, N and AF are empty placeholder cells.
$extraValueCols = @("K","AC")
$extraBlankCols = @("J","L","N","AF")

function Get-RowValues($row) {
    $vals = @{}
    foreach ($col in $dataCols) {
        $vals[$col] = $ws.Range("$col$row").Value()
    }
    return $vals
}

# Y (Startdatum) and AA (Slutdatum) hold plain date-looking text (e.g.
# "2023-08-15"); Excel auto-converts such text to a real date serial when it
# is assigned through .Value. Force the destination cell to a text format
# first so the original text is preserved verbatim.
$dateTextCols = @("Y","AA")

function Set-RowValues($row, $vals) {
    foreach ($col in $dataCols) {
        $cell = $ws.Range("$col$row")
        if ($dateTextCols -contains $col) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $vals[$col]
    }
}

# First move the empty placeholder cells from row 6 to row 4 while the rows
# still hold their original data (Cut preserves the "present but empty" cell
# instead of just deleting/recreating it).
foreach ($col in $extraBlankCols) {
    $ws.Range("$col`6").Cut($ws.Range("$col`4")) | Out-Null
}

# Capture the remaining (non-blank) values of rows 2, 3, 4 and 6 before any
# of them get overwritten.
$row2 = Get-RowValues 2
$row3 = Get-RowValues 3
$row4 = Get-RowValues 4
$row6 = Get-RowValues 6

$row6Extra = @{}
foreach ($col in $extraValueCols) {
    $row6Extra[$col] = $ws.Range("$col`6").Value()
}

# Apply the rotation.
Set-RowValues 2 $row4
Set-RowValues 3 $row2
Set-RowValues 4 $row6
Set-RowValues 6 $row3

# Row 4 now also carries the K/AC values that used to belong to row 6.
foreach ($col in $extraValueCols) {
    $ws.Range("$col`4").Value = $row6Extra[$col]
}

# Row 6 no longer carries the K/AC values (they now live on row 4).
foreach ($col in $extraValueCols) {
    $ws.Range("$col`6").ClearContents() | Out-Null
}
